$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 174; existing rows 174..351 shift down to 175..352,
# preserving all of their original data/formatting.
$ws.Rows.Item(174).Insert()

# Populate the newly inserted row 174 with the new weekly record.
$ws.Cells.Item(174, 1).Value = 3
$ws.Cells.Item(174, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(174, 3).Value = 'Coquimbo'
$ws.Cells.Item(174, 4).Value = 44789
$ws.Cells.Item(174, 5).Value = 5
$ws.Cells.Item(174, 6).Value = 100112039
$ws.Cells.Item(174, 7).Value = 'Ciboulette'
$ws.Cells.Item(174, 8).Value = 'Sin especificar'
$ws.Cells.Item(174, 9).Value = 'Primera'
$ws.Cells.Item(174, 10).Value = 120
$ws.Cells.Item(174, 11).Value = 1500
$ws.Cells.Item(174, 12).Value = 1500
$ws.Cells.Item(174, 13).Value = 1500
$ws.Cells.Item(174, 14).Value = '$/docena de atados'
$ws.Cells.Item(174, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(174, 16).Value = 500
$ws.Cells.Item(174, 17).Value = 3
$ws.Cells.Item(174, 18).Value = 'Hortaliza'
